# Rendering Menu and cart components
#
# This script reproduces the authored edit against "Aaswad Caterers.xlsx":
#  - Tasks sheet: new backlog items appended (rows 52-55, 57)
#  - Schema Design sheet: "orders" -> "orders[]" in the Cart Items Schema
#  - View state: Schema Design becomes the active/selected tab, with the
#    Tasks and Schema Design sheets' selections/scroll position updated.

$wb = $excel.ActiveWorkbook

$wsTasks  = $wb.Worksheets.Item("Tasks")
$wsLogs   = $wb.Worksheets.Item("Logs")
$wsSchema = $wb.Worksheets.Item("Schema Design")

# --- Tasks sheet: new Backlog entries -------------------------------------
$wsTasks.Range("A52").Value = "price range"
$wsTasks.Range("A53").Value = "sort order up and down arrow "
$wsTasks.Range("A54").Value = "category add selection dropdown"
$wsTasks.Range("A55").Value = "in view add update button and delete button"
$wsTasks.Range("A57").Value = "pagination for items"

# --- Schema Design sheet: Cart Items Schema tweak --------------------------
$wsSchema.Range("B14").Value = "orders[]"

# --- View / selection state -------------------------------------------------
# Tasks: scroll down to the newly added rows and select B53.
$wsTasks.Activate()
$wsTasks.Range("B53").Select()

# Logs: no longer the selected tab; keep its existing selection (C13).
$wsLogs.Activate()

# Schema Design: becomes the active/selected tab, selection moves to C15.
$wsSchema.Activate()
$wsSchema.Range("C15").Select()
